# contraordenarcheque.xlsx data refactor
# Updates the cheque-range test values on sheet "Datos" and restores the
# sheet's scroll position, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Data changes (numeroCheque / rangoDesde / rangoHasta columns) ---
$ws.Range("P2").Value = 65765
$ws.Range("Q3").Value = 65762
$ws.Range("R3").Value = 65763

# --- View changes: scroll one column left (was topLeftCell J1, now I1) ---
$ws.Activate()
$ws.Range("I1").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1

# Restore the original active cell/selection (P2), as in the source file.
$ws.Range("P2").Select()
